$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Report header: new weekly volume/number and the Week-Covering date range.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# ---------------------------------------------------------------------------
# Crime-complaint table refresh (rows 14-29, columns C:N).
#
# Some cells flip between a numeric figure and the sheet's "no data"
# placeholder text ("0" / "***.*"), which also carries its own cell style.
# Rather than just writing .Value (which would coerce a numeric-looking
# string back into a Number and keep the old style), those transitions are
# done by copying a same-styled placeholder/numeric cell that is never
# itself edited by this script (C14/E14 hold the literal "0"/"***.*" text
# placeholders; I14/K14 hold plain numeric styles) onto the target cell.
# ---------------------------------------------------------------------------

# Row 14 - Murder: Week-to-Date 28-day comparison no longer meaningful -> placeholders
$ws.Range("C14").Copy($ws.Range("G14"))   # -> "0"
$ws.Range("E14").Copy($ws.Range("H14"))   # -> "***.*"

# Row 15 - Rape
$ws.Range("M15").Value = 70

# Row 16 - Robbery
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -54.545454545454
$ws.Range("I16").Value = 170
$ws.Range("J16").Value = 108
$ws.Range("K16").Value = 57.407407407407
$ws.Range("L16").Value = 61.904761904761
$ws.Range("M16").Value = 32.8125
$ws.Range("N16").Value = -78.288633461047

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = -58.333333333333
$ws.Range("I17").Value = 198
$ws.Range("J17").Value = 191
$ws.Range("K17").Value = 3.664921465968
$ws.Range("L17").Value = 69.230769230769
$ws.Range("M17").Value = 62.295081967213
$ws.Range("N17").Value = -8.755760368663

# Row 18 - Burglary
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -65
$ws.Range("I18").Value = 152
$ws.Range("J18").Value = 118
$ws.Range("K18").Value = 28.813559322033
$ws.Range("L18").Value = -9.523809523809
$ws.Range("M18").Value = 58.333333333333
$ws.Range("N18").Value = -58.918918918918

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -56.521739130434
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = -26.760563380281
$ws.Range("I19").Value = 730
$ws.Range("J19").Value = 553
$ws.Range("K19").Value = 32.007233273056
$ws.Range("L19").Value = 70.560747663551
$ws.Range("M19").Value = 169.372693726937
$ws.Range("N19").Value = 67.048054919908

# Row 20 - G.L.A.: previously-empty columns now have figures
$ws.Range("I14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 2
$ws.Range("I14").Copy($ws.Range("F20"))
$ws.Range("F20").Value = 2
$ws.Range("I20").Value = 58
$ws.Range("K20").Value = 11.538461538461
$ws.Range("L20").Value = 7.407407407407
$ws.Range("M20").Value = 23.404255319148
$ws.Range("N20").Value = -84.696569920844

# Row 21 - TOTAL
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -38.235294117647
$ws.Range("G21").Value = 127
$ws.Range("H21").Value = -39.370078740157
$ws.Range("I21").Value = 1327
$ws.Range("J21").Value = 1036
$ws.Range("K21").Value = 28.088803088803
$ws.Range("L21").Value = 49.943502824858
$ws.Range("M21").Value = 96.592592592592
$ws.Range("N21").Value = -40.090293453724

# Row 22 - Transit
$ws.Range("I14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("C14").Copy($ws.Range("F22"))   # -> "0"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 50

# Row 23 - Housing
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 163
$ws.Range("J23").Value = 156
$ws.Range("K23").Value = 4.487179487179
$ws.Range("L23").Value = -3.550295857988
$ws.Range("M23").Value = 28.346456692913

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = 6.976744186046
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 163
$ws.Range("H24").Value = -28.220858895705
$ws.Range("I24").Value = 2119
$ws.Range("J24").Value = 1525
$ws.Range("K24").Value = 38.950819672131
$ws.Range("L24").Value = 174.837872892348
$ws.Range("M24").Value = 197.194950911641

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -6.25
$ws.Range("I25").Value = 429
$ws.Range("J25").Value = 316
$ws.Range("K25").Value = 35.759493670886
$ws.Range("L25").Value = 25.072886297376
$ws.Range("M25").Value = 36.624203821656

# Row 26 - UCR Rape*
$ws.Range("I14").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1
$ws.Range("C14").Copy($ws.Range("D26"))   # -> "0"
$ws.Range("E14").Copy($ws.Range("E26"))   # -> "***.*"
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 27
$ws.Range("K26").Value = 80
$ws.Range("L26").Value = 28.571428571428

# Row 27 - Other Sex Crimes
$ws.Range("C14").Copy($ws.Range("C27"))   # -> "0"
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = 70.370370370370

# Row 28 - Shooting Vic.
$ws.Range("C14").Copy($ws.Range("G28"))   # -> "0"
$ws.Range("E14").Copy($ws.Range("H28"))   # -> "***.*"

# Row 29 - Shooting Inc.
$ws.Range("C14").Copy($ws.Range("G29"))   # -> "0"
$ws.Range("E14").Copy($ws.Range("H29"))   # -> "***.*"
